$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("targets")

# Insert a new row above row 2, shifting existing rows 2-13 down to 3-14.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the new cohort entry.
$ws.Cells.Item(2, 1).Value = 878
$ws.Cells.Item(2, 2).Value = "Non-Emergent MNCS (age 18 or greater), post op Afib (any)"

# Insert() copies the formatting of the row above (the bold/centered header),
# so reset the new data row back to the plain "Normal" style used by the rest
# of the data rows.
$ws.Range("A2:B2").Style = "Normal"

# The old last row (originally row 13: id 1106, "Non-Emergent Major Non Cardiac
# Surgery no prior Opioid") has shifted to row 14 and is no longer part of the
# data set, so remove it entirely to restore the original row count (A1:B13).
$ws.Rows.Item(14).Delete()
